# Regenerate the lattice-multiplication exercise table: each of the 15
# cells gets a brand new "A x B" problem (title line, the two digits of
# B spaced out, a divider, and two placeholder rows seeded with the
# digits of A). The table shape (5 rows x 3 cols) is unchanged - only
# the text inside each cell's single run changes.

$d = $word.ActiveDocument
$t = $d.Tables(1)
$vt = [char]11   # vertical-tab char -> becomes <w:br/> between the w:t runs

# New content per cell, in row-major order (row, col, 5 lines of text).
$cells = @(
    @(1,1,"25 x 71","  7    1","  ----","2|    |","5|    |"),
    @(1,2,"38 x 34","  3    4","  ----","3|    |","8|    |"),
    @(1,3,"48 x 44","  4    4","  ----","4|    |","8|    |"),
    @(2,1,"73 x 76","  7    6","  ----","7|    |","3|    |"),
    @(2,2,"20 x 92","  9    2","  ----","2|    |","0|    |"),
    @(2,3,"89 x 75","  7    5","  ----","8|    |","9|    |"),
    @(3,1,"20 x 55","  5    5","  ----","2|    |","0|    |"),
    @(3,2,"73 x 68","  6    8","  ----","7|    |","3|    |"),
    @(3,3,"10 x 31","  3    1","  ----","1|    |","0|    |"),
    @(4,1,"37 x 19","  1    9","  ----","3|    |","7|    |"),
    @(4,2,"34 x 81","  8    1","  ----","3|    |","4|    |"),
    @(4,3,"95 x 90","  9    0","  ----","9|    |","5|    |"),
    @(5,1,"73 x 68","  6    8","  ----","7|    |","3|    |"),
    @(5,2,"78 x 27","  2    7","  ----","7|    |","8|    |"),
    @(5,3,"45 x 90","  9    0","  ----","4|    |","5|    |")
)

foreach ($c in $cells) {
    $row = $c[0]
    $col = $c[1]
    $lines = $c[2..6]
    $newText = [string]::Join($vt, $lines)
    $t.Cell($row, $col).Range.Text = $newText
}

Write-Output "Updated $($cells.Count) cells"
